# "big data dan cloud computing" commit: PowerPoint's cached
# datetimeFigureOut placeholders (slide master + every slide layout)
# get re-stamped from the old capture date (10/7/24) to the new one
# (11/8/24), and the deck picks up an (empty) slide-guide extension
# list that PowerPoint writes into presentation.xml on save.

$p = $ppt.ActivePresentation

$newDate = "11/8/24"

$m = $p.SlideMaster

# Update the "Date Placeholder" field on the slide master itself.
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

# Update the "Date Placeholder" field on every slide layout belonging
# to the master (there are 11 of them in this deck).
for ($j = 1; $j -le $m.CustomLayouts.Count; $j++) {
    $layout = $m.CustomLayouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}
